# Auto-generated edit script applying the 2026-01-17 06:45:50 scraper refresh
# for "LÍNEA 141" horarios workbook (sheets LP1912, LP1912-215, 6203-6173).
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: LP1912
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("A2").Value = "Última actualización: 06:45:50"
$ws1.Range("A3").Value = "Total filas: 47"

# 6 new rows appear in the refreshed scrape (41 -> 47 data rows); insert blank
# rows at the first changed position so everything below shifts down, matching
# the new layout, then rewrite rows 19-52 with the refreshed values.
$ws1.Range("A19:A24").EntireRow.Insert()

$ws1.Cells.Item(19, 1).Value = "06:45:50"
$ws1.Cells.Item(19, 2).Value = "07:00"
$ws1.Cells.Item(19, 3).Value = "16_SANTA ANA"
$ws1.Cells.Item(19, 4).Value = 15
$ws1.Cells.Item(19, 5).Value = "LP1912"
$ws1.Cells.Item(20, 1).Value = "05:57:08"
$ws1.Cells.Item(20, 2).Value = "07:01"
$ws1.Cells.Item(20, 3).Value = "16_SANTA ANA"
$ws1.Cells.Item(20, 4).Value = 64
$ws1.Cells.Item(20, 5).Value = "LP1912"
$ws1.Cells.Item(21, 1).Value = "06:45:50"
$ws1.Cells.Item(21, 2).Value = "07:03"
$ws1.Cells.Item(21, 3).Value = "225_GOMEZ"
$ws1.Cells.Item(21, 4).Value = 18
$ws1.Cells.Item(21, 5).Value = "LP1912"
$ws1.Cells.Item(22, 1).Value = "05:42:22"
$ws1.Cells.Item(22, 2).Value = "07:04"
$ws1.Cells.Item(22, 3).Value = "225_GOMEZ"
$ws1.Cells.Item(22, 4).Value = 82
$ws1.Cells.Item(22, 5).Value = "LP1912"
$ws1.Cells.Item(23, 1).Value = "06:33:46"
$ws1.Cells.Item(23, 2).Value = "07:06"
$ws1.Cells.Item(23, 3).Value = "215C_EL PATO"
$ws1.Cells.Item(23, 4).Value = 33
$ws1.Cells.Item(23, 5).Value = "LP1912"
$ws1.Cells.Item(24, 1).Value = "05:42:22"
$ws1.Cells.Item(24, 2).Value = "07:07"
$ws1.Cells.Item(24, 3).Value = "215C_EL PATO"
$ws1.Cells.Item(24, 4).Value = 85
$ws1.Cells.Item(24, 5).Value = "LP1912"
$ws1.Cells.Item(25, 1).Value = "06:33:46"
$ws1.Cells.Item(25, 2).Value = "07:13"
$ws1.Cells.Item(25, 3).Value = "14X44_ABASTO"
$ws1.Cells.Item(25, 4).Value = 40
$ws1.Cells.Item(25, 5).Value = "LP1912"
$ws1.Cells.Item(26, 1).Value = "05:42:22"
$ws1.Cells.Item(26, 2).Value = "07:14"
$ws1.Cells.Item(26, 3).Value = "14X44_ABASTO"
$ws1.Cells.Item(26, 4).Value = 92
$ws1.Cells.Item(26, 5).Value = "LP1912"
$ws1.Cells.Item(27, 1).Value = "06:33:46"
$ws1.Cells.Item(27, 2).Value = "07:15"
$ws1.Cells.Item(27, 3).Value = "16_SANTA ANA"
$ws1.Cells.Item(27, 4).Value = 42
$ws1.Cells.Item(27, 5).Value = "LP1912"
$ws1.Cells.Item(28, 1).Value = "06:33:46"
$ws1.Cells.Item(28, 2).Value = "07:20"
$ws1.Cells.Item(28, 3).Value = "215A_EL PATO"
$ws1.Cells.Item(28, 4).Value = 47
$ws1.Cells.Item(28, 5).Value = "LP1912"
$ws1.Cells.Item(29, 1).Value = "05:42:22"
$ws1.Cells.Item(29, 2).Value = "07:21"
$ws1.Cells.Item(29, 3).Value = "215A_EL PATO"
$ws1.Cells.Item(29, 4).Value = 99
$ws1.Cells.Item(29, 5).Value = "LP1912"
$ws1.Cells.Item(30, 1).Value = "05:57:08"
$ws1.Cells.Item(30, 2).Value = "07:29"
$ws1.Cells.Item(30, 3).Value = "14_ABASTO"
$ws1.Cells.Item(30, 4).Value = 92
$ws1.Cells.Item(30, 5).Value = "LP1912"
$ws1.Cells.Item(31, 1).Value = "05:42:22"
$ws1.Cells.Item(31, 2).Value = "07:33"
$ws1.Cells.Item(31, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(31, 4).Value = 111
$ws1.Cells.Item(31, 5).Value = "LP1912"
$ws1.Cells.Item(32, 1).Value = "05:57:08"
$ws1.Cells.Item(32, 2).Value = "07:34"
$ws1.Cells.Item(32, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(32, 4).Value = 97
$ws1.Cells.Item(32, 5).Value = "LP1912"
$ws1.Cells.Item(33, 1).Value = "06:45:50"
$ws1.Cells.Item(33, 2).Value = "07:35"
$ws1.Cells.Item(33, 3).Value = "17X38_ROMERO"
$ws1.Cells.Item(33, 4).Value = 50
$ws1.Cells.Item(33, 5).Value = "LP1912"
$ws1.Cells.Item(34, 1).Value = "05:42:22"
$ws1.Cells.Item(34, 2).Value = "07:36"
$ws1.Cells.Item(34, 3).Value = "17X38_ROMERO"
$ws1.Cells.Item(34, 4).Value = 114
$ws1.Cells.Item(34, 5).Value = "LP1912"
$ws1.Cells.Item(35, 1).Value = "06:33:46"
$ws1.Cells.Item(35, 2).Value = "07:36"
$ws1.Cells.Item(35, 3).Value = "27_EL RETIRO"
$ws1.Cells.Item(35, 4).Value = 63
$ws1.Cells.Item(35, 5).Value = "LP1912"
$ws1.Cells.Item(36, 1).Value = "05:42:22"
$ws1.Cells.Item(36, 2).Value = "07:37"
$ws1.Cells.Item(36, 3).Value = "27_EL RETIRO"
$ws1.Cells.Item(36, 4).Value = 115
$ws1.Cells.Item(36, 5).Value = "LP1912"
$ws1.Cells.Item(37, 1).Value = "06:33:46"
$ws1.Cells.Item(37, 2).Value = "07:43"
$ws1.Cells.Item(37, 3).Value = "10_OLMOS"
$ws1.Cells.Item(37, 4).Value = 70
$ws1.Cells.Item(37, 5).Value = "LP1912"
$ws1.Cells.Item(38, 1).Value = "05:57:08"
$ws1.Cells.Item(38, 2).Value = "07:44"
$ws1.Cells.Item(38, 3).Value = "10_OLMOS"
$ws1.Cells.Item(38, 4).Value = 107
$ws1.Cells.Item(38, 5).Value = "LP1912"
$ws1.Cells.Item(39, 1).Value = "05:57:08"
$ws1.Cells.Item(39, 2).Value = "07:51"
$ws1.Cells.Item(39, 3).Value = "15_ABASTO"
$ws1.Cells.Item(39, 4).Value = 114
$ws1.Cells.Item(39, 5).Value = "LP1912"
$ws1.Cells.Item(40, 1).Value = "06:16:15"
$ws1.Cells.Item(40, 2).Value = "07:58"
$ws1.Cells.Item(40, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(40, 4).Value = 102
$ws1.Cells.Item(40, 5).Value = "LP1912"
$ws1.Cells.Item(41, 1).Value = "06:33:46"
$ws1.Cells.Item(41, 2).Value = "07:59"
$ws1.Cells.Item(41, 3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(41, 4).Value = 86
$ws1.Cells.Item(41, 5).Value = "LP1912"
$ws1.Cells.Item(42, 1).Value = "06:16:15"
$ws1.Cells.Item(42, 2).Value = "08:00"
$ws1.Cells.Item(42, 3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(42, 4).Value = 104
$ws1.Cells.Item(42, 5).Value = "LP1912"
$ws1.Cells.Item(43, 1).Value = "06:33:46"
$ws1.Cells.Item(43, 2).Value = "08:00"
$ws1.Cells.Item(43, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(43, 4).Value = 87
$ws1.Cells.Item(43, 5).Value = "LP1912"
$ws1.Cells.Item(44, 1).Value = "06:45:50"
$ws1.Cells.Item(44, 2).Value = "08:01"
$ws1.Cells.Item(44, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(44, 4).Value = 76
$ws1.Cells.Item(44, 5).Value = "LP1912"
$ws1.Cells.Item(45, 1).Value = "06:16:15"
$ws1.Cells.Item(45, 2).Value = "08:03"
$ws1.Cells.Item(45, 3).Value = "17_ROMERO"
$ws1.Cells.Item(45, 4).Value = 107
$ws1.Cells.Item(45, 5).Value = "LP1912"
$ws1.Cells.Item(46, 1).Value = "06:33:46"
$ws1.Cells.Item(46, 2).Value = "08:12"
$ws1.Cells.Item(46, 3).Value = "10_OLMOS"
$ws1.Cells.Item(46, 4).Value = 99
$ws1.Cells.Item(46, 5).Value = "LP1912"
$ws1.Cells.Item(47, 1).Value = "06:16:15"
$ws1.Cells.Item(47, 2).Value = "08:15"
$ws1.Cells.Item(47, 3).Value = "17_ROMERO"
$ws1.Cells.Item(47, 4).Value = 119
$ws1.Cells.Item(47, 5).Value = "LP1912"
$ws1.Cells.Item(48, 1).Value = "06:33:46"
$ws1.Cells.Item(48, 2).Value = "08:26"
$ws1.Cells.Item(48, 3).Value = "15X38_ABASTO"
$ws1.Cells.Item(48, 4).Value = 113
$ws1.Cells.Item(48, 5).Value = "LP1912"
$ws1.Cells.Item(49, 1).Value = "06:33:46"
$ws1.Cells.Item(49, 2).Value = "08:27"
$ws1.Cells.Item(49, 3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws1.Cells.Item(49, 4).Value = 114
$ws1.Cells.Item(49, 5).Value = "LP1912"
$ws1.Cells.Item(50, 1).Value = "06:45:50"
$ws1.Cells.Item(50, 2).Value = "08:29"
$ws1.Cells.Item(50, 3).Value = "14_ABASTO"
$ws1.Cells.Item(50, 4).Value = 104
$ws1.Cells.Item(50, 5).Value = "LP1912"
$ws1.Cells.Item(51, 1).Value = "06:33:46"
$ws1.Cells.Item(51, 2).Value = "08:31"
$ws1.Cells.Item(51, 3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(51, 4).Value = 118
$ws1.Cells.Item(51, 5).Value = "LP1912"
$ws1.Cells.Item(52, 1).Value = "06:45:50"
$ws1.Cells.Item(52, 2).Value = "08:38"
$ws1.Cells.Item(52, 3).Value = "215C_EL PATO"
$ws1.Cells.Item(52, 4).Value = 113
$ws1.Cells.Item(52, 5).Value = "LP1912"

# ---------------------------------------------------------------------
# Sheet 2: LP1912-215
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("A2").Value = "Última actualización: 06:45:50"
$ws2.Range("A3").Value = "Total filas: 8"

# New row 13 appended with the refreshed 215C_EL PATO entry
$ws2.Cells.Item(13, 1).Value = "06:45:50"
$ws2.Cells.Item(13, 2).Value = "08:38"
$ws2.Cells.Item(13, 3).Value = "215C_EL PATO"
$ws2.Cells.Item(13, 4).Value = 113
$ws2.Cells.Item(13, 5).Value = "LP1912"

# ---------------------------------------------------------------------
# Sheet 3: 6203-6173
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Range("A2").Value = "Última actualización: 06:45:50"
